$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename strategy labels in column B (rows 2-9) ---
$ws.Range("B2").Value = "wiley"
$ws.Range("B4").Value = "springer"
$ws.Range("B5").Value = "scopus"
$ws.Range("B6").Value = "sciencedirect"
$ws.Range("B7").Value = "ieee"
$ws.Range("B8").Value = "googlescholar"
$ws.Range("B9").Value = "acm"
# B3 (webofscience) stays the same

# --- Row 2 (wiley) ---
$ws.Range("C2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("K2").Value = 15
$ws.Range("N2").Value = 15
$ws.Range("P2").Value = 15

# --- Row 4 (springer) ---
$ws.Range("C4").Value = 1.54
$ws.Range("D4").Value = 7.140000000000001
$ws.Range("E4").Value = 2.53
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 65
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 65
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 65
$ws.Range("Q4").Value = 1
$ws.Range("S4").Value = 1.54
$ws.Range("T4").Value = 7.140000000000001
$ws.Range("U4").Value = 2.53
$ws.Range("V4").Value = 0.01538461538461539
$ws.Range("W4").Value = 0.07142857142857142
$ws.Range("X4").Value = 0.02531645569620253

# --- Row 5 (scopus) ---
$ws.Range("C5").Value = 3.8
$ws.Range("D5").Value = 21.43
$ws.Range("E5").Value = 6.45
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 79
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = 79
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 2
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 77
$ws.Range("S5").Value = 3.8
$ws.Range("T5").Value = 21.43
$ws.Range("U5").Value = 6.45
$ws.Range("V5").Value = 0.0379746835443038
$ws.Range("W5").Value = 0.2142857142857143
$ws.Range("X5").Value = 0.06451612903225806

# --- Row 6 (sciencedirect) ---
$ws.Range("C6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("K6").Value = 21
$ws.Range("N6").Value = 21
$ws.Range("P6").Value = 21
$ws.Range("S6").Value = 0
$ws.Range("V6").Value = 0

# --- Row 7 (ieee) ---
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").ClearContents()
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("S7").Value = 2.36
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0.02360515021459228
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0

# --- Row 8 (googlescholar) ---
$ws.Range("C8").Value = 2.36
$ws.Range("D8").Value = 78.56999999999999
$ws.Range("E8").Value = 4.58
$ws.Range("J8").Value = 11
$ws.Range("K8").Value = 466
$ws.Range("M8").Value = 11
$ws.Range("N8").Value = 466
$ws.Range("O8").Value = 11
$ws.Range("P8").Value = 466
$ws.Range("Q8").Value = 11
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 2.36
$ws.Range("T8").Value = 78.56999999999999
$ws.Range("U8").Value = 4.58
$ws.Range("V8").Value = 0.02360515021459228
$ws.Range("W8").Value = 0.7857142857142857
$ws.Range("X8").Value = 0.04583333333333334

# --- Row 9 (acm) ---
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").ClearContents()
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 2.29
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0.02291666666666667
$ws.Range("W9").Value = 0
$ws.Range("X9").Value = 0

# --- Row 10 (s0 / union) ---
$ws.Range("S10").Value = 2.29
$ws.Range("T10").Value = 78.56999999999999
$ws.Range("U10").Value = 4.45
